# VACCINECERT-1633 Fixed CSV upload templates
#
# The CSV "post" template for vaccination/antibody uploads had a stale
# sample "sampleDate" (I2). Update it to the corrected date and leave the
# cursor where the author left it (J7) when they saved the file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# I2 held the sample sampleDate as a date serial (2021-11-01 = 44501).
# Fix it to the corrected date (2021-11-16 = 44516).
$ws.Range("I2").Value = 44516

# Reflect the author's final cursor position on the sheet.
$ws.Range("J7").Select()
